$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "68.929.86"
$ws.Cells.Item(2, 5).Value = "  -3.55%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.508.44"
$ws.Cells.Item(3, 5).Value = "  -3.44%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.22%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "575.26"
$ws.Cells.Item(5, 5).Value = "  -2.22%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "170.36"
$ws.Cells.Item(6, 5).Value = "  -5.39%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.614"
$ws.Cells.Item(7, 5).Value = "  -0.07%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "3.505.78"
$ws.Cells.Item(8, 5).Value = "  -3.26%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  +0.02%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.189"
$ws.Cells.Item(10, 5).Value = "  -6.57%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "6.53"
$ws.Cells.Item(11, 5).Value = "  +11.31%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.597"
$ws.Cells.Item(12, 5).Value = "  -1.43%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "47.04"
$ws.Cells.Item(13, 5).Value = "  -5.09%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.0000274"
$ws.Cells.Item(14, 5).Value = "  -4.10%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "687.29"
$ws.Cells.Item(15, 5).Value = "  +0.78%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "4.079.06"
$ws.Cells.Item(16, 5).Value = "  -3.31%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "8.74"
$ws.Cells.Item(17, 5).Value = "  -2.65%  "

# Row 18
$ws.Cells.Item(18, 2).Value = "WrappedEther"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(18, 4).Value = "3.534.88"
$ws.Cells.Item(18, 5).Value = "  -2.66%  "

# Row 19
$ws.Cells.Item(19, 2).Value = "WrappedBTC"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(19, 4).Value = "68.984.59"
$ws.Cells.Item(19, 5).Value = "  -3.56%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  -1.46%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "17.31"
$ws.Cells.Item(21, 5).Value = "  -5.12%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "11.13"
$ws.Cells.Item(22, 5).Value = "  -4.23%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.909"
$ws.Cells.Item(23, 5).Value = "  -2.58%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "16.49"
$ws.Cells.Item(24, 5).Value = "  -7.12%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "97.30"
$ws.Cells.Item(25, 5).Value = "  -5.53%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "3.81"
$ws.Cells.Item(26, 5).Value = "  -4.76%  "

# Row 27
$ws.Cells.Item(27, 2).Value = "Dai"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.999"
$ws.Cells.Item(27, 5).Value = "  -0.17%  "

# Row 28
$ws.Cells.Item(28, 2).Value = "ImmutableX"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "2.67"
$ws.Cells.Item(28, 5).Value = "  -5.61%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "9.39"
$ws.Cells.Item(29, 5).Value = "  -5.63%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "33.26"
$ws.Cells.Item(30, 5).Value = "  -4.86%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "8.82"
$ws.Cells.Item(31, 5).Value = "  -3.95%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.16"
$ws.Cells.Item(32, 5).Value = "  -6.97%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "7.24"
$ws.Cells.Item(33, 5).Value = "  -0.11%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  -5.25%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "3.75"
$ws.Cells.Item(35, 5).Value = "  -8.95%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "568.24"
$ws.Cells.Item(36, 5).Value = "  -1.07%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "10.82"
$ws.Cells.Item(37, 5).Value = "  -4.31%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  -4.25%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "57.36"
$ws.Cells.Item(39, 5).Value = "  -3.56%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  +0.36%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "3.481.71"
$ws.Cells.Item(41, 5).Value = "  -5.03%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.0438"
$ws.Cells.Item(42, 5).Value = "  -6.29%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "TheGraph"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.336"
$ws.Cells.Item(43, 5).Value = "  -2.78%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "Kaspa"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.136"
$ws.Cells.Item(44, 5).Value = "  -4.24%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "33.28"
$ws.Cells.Item(45, 5).Value = "  -6.37%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "0.0₃0703"
$ws.Cells.Item(46, 5).Value = "  -7.29%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  +3.94%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "2.57"
$ws.Cells.Item(48, 5).Value = "  -6.39%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.132"
$ws.Cells.Item(49, 5).Value = "  -1.13%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "133.62"
$ws.Cells.Item(50, 5).Value = "  +1.71%  "

# Row 51
$ws.Cells.Item(51, 5).Value = "  -0.76%  "

